{"js": "// Chuyen de tai nhom tu \"App viet nhat ky\" sang \"Game Color Bubble\".\n// \"Bubble\" keeps its own distinct run formatting (Segoe UI, color\n// #24292E, white run shading) exactly like the source diff - it is a\n// second run appended right after \"Game Color \".\n\n// The paragraph's trailing \"_GoBack\" bookmark sits immediately after the\n// run we're replacing. Remove it first so it doesn't end up wrapped\n// around the freshly inserted runs; we recreate it (via the OOXML\n// payload below) right after them, matching the original structure.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n} catch (e) {\n  // no existing _GoBack bookmark - nothing to delete\n}\nawait context.sync();\n\n// Locate the exact run to replace.\nconst results = context.document.body.search(\"App vi\u1ebft nh\u1eadt k\u00fd\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"App vi\u1ebft nh\u1eadt k\u00fd\" in the document.');\n}\n\nconst target = results.items[0];\n\n// Replace that single run with two runs:\n//   1) \"Game Color \" - keeps the original TimesNewRomanPSMT / sz24 formatting\n//   2) \"Bubble\"       - new Segoe UI / #24292E / white-shaded run\n// followed by a fresh _GoBack bookmark, mirroring the original layout.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r w:rsidR=\"00BB3C18\">' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"TimesNewRomanPSMT\" w:hAnsi=\"TimesNewRomanPSMT\" w:cs=\"TimesNewRomanPSMT\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\">Game Color </w:t>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n  '<w:color w:val=\"24292E\"/>' +\n  '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>' +\n  '</w:rPr>' +\n  '<w:t>Bubble</w:t>' +\n  '</w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Chuyen de tai nhom tu \"App viet nhat ky\" sang \"Game Color Bubble\".\n# The word \"Bubble\" keeps its own (GitHub-ish) run formatting: Segoe UI,\n# color #24292E and a white run shading - exactly like the source diff.\n\n$d = $word.ActiveDocument\n\n# The trailing \"_GoBack\" bookmark sits right after the target run; remove\n# it first so the replacement lands *before* a freshly created one instead\n# of the new runs getting sandwiched between bookmarkStart/bookmarkEnd.\ntry {\n    $goBack = $d.Bookmarks.Item(\"_GoBack\")\n    $goBack.Delete()\n} catch {\n    # no existing _GoBack bookmark - nothing to do\n}\n\n# Locate the exact run we need to replace.\n$findRange = $d.Content\n$findRange.Find.Execute(\"App vi\u1ebft nh\u1eadt k\u00fd\") | Out-Null\n$start = $findRange.Start\n$end = $findRange.End\n\n$target = $d.Range($start, $end)\n\n# Replace the single \"App vi\u1ebft nh\u1eadt k\u00fd\" run with two runs:\n#   1) \"Game Color \" - keeps the original TimesNewRomanPSMT/sz24 formatting\n#   2) \"Bubble\"       - new Segoe UI / #24292E / white-shaded run\n# then restore the _GoBack bookmark right after them, matching the\n# original document's structure.\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidR=\"00BB3C18\">\n              <w:rPr>\n                <w:rFonts w:ascii=\"TimesNewRomanPSMT\" w:hAnsi=\"TimesNewRomanPSMT\" w:cs=\"TimesNewRomanPSMT\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Game Color </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>\n                <w:color w:val=\"24292E\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n              </w:rPr>\n              <w:t>Bubble</w:t>\n            </w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($xml)\n"}
